$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (23) had the "last row" date style (plain date).
# Now that a new row is appended, row 23 reverts to the regular date+time
# style used by all the other interior rows, and the new row 24 becomes
# the new "last row" with the plain date style.
$ws.Cells.Item(23, 1).NumberFormat = $ws.Cells.Item(22, 1).NumberFormat

# Append the new day's data (row 24).
$ws.Cells.Item(24, 1).Value = 45608
$ws.Cells.Item(24, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(24, 2).Value = 59
$ws.Cells.Item(24, 3).Value = 50
$ws.Cells.Item(24, 4).Value = 55
